$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the RF column (I) for rows 17 through 39 from 37.25942528735632 to 51.2
$ws.Range("I17:I39").Value = 51.2
